# Updated User stories excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: cell values below are written in the same left-to-right, row-by-row
# order the original author used so the shared-string table is built up in
# the same sequence as the authored workbook.

# --- Row 6 (Owen / Create Forum Topic) new detail columns ---
$ws.Range("D6").Value = "Create a new topic"
$ws.Range("E6").Value = "To provide a forum for users to post in"
$ws.Range("F6").Value = "There is not an identical forum topic"
$ws.Range("G6").Value = "Users can post content inside of the topic"

# --- Row 7 (Owen / View Content) new detail columns ---
$ws.Range("E7").Value = "To see content that other users have posted"
$ws.Range("F7").Value = "Content has been posted"
$ws.Range("G7").Value = "Users can respond to post, Administrator can edit and delete posts"

# --- Row 8 (Owen / Create Account) new detail columns ---
$ws.Range("D8").Value = "Allows user to create their own personal account"
$ws.Range("E8").Value = "To allow users to join and contribute to the community"
$ws.Range("F8").Value = "User does not have an account"
$ws.Range("G8").Value = "User now has access to privilege given to their user class"

$ws.Range("C7").Value = "Administrator, Content Publisher, Registered User, Unregistered User"
$ws.Range("C8").Value = "Unregistered User"
$ws.Range("D7").Value = "Allow all users to view content posted"

# --- Row 6 & 7 & 8 column C (Administrator) is a reused string from elsewhere ---
$ws.Range("C6").Value = "Administrator"

# --- New column H: "Possible Breakdown?" ---
$ws.Range("H1").Value = "Possible Breakdown?"
$ws.Range("H3").Value = "Yes"
$ws.Range("H4").Value = "Yes"

# Wrap text on the newly populated detail cells (matches existing style used
# by the other rows in this table).
$ws.Range("C6:G8").WrapText = $true

# Row heights for the newly-expanded rows
$ws.Rows(6).RowHeight = 45
$ws.Rows(7).RowHeight = 75
$ws.Rows(8).RowHeight = 75

# --- Column widths ---
$ws.Columns(3).ColumnWidth = 24
$ws.Columns(8).ColumnWidth = 21

# --- Sheet view: zoom + selection ---
$ws.Select()
$excel.ActiveWindow.Zoom = 70
$ws.Range("H2").Select()
